# Add two new columns "I0" (column I) and "IF" (column J) to the sheet,
# matching the style of the existing header row and filling in the
# per-row numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1) ---
# First put the text in, then clone the formatting of the existing
# header cell H1 (bold font, border, centered alignment) onto I1:J1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Re-apply the values in case PasteSpecial touched the cell contents.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data cells (rows 2-24) ---
$iValues = @(1,1,5,6,7,6,6,6,3,7,6,6,8,4,5,8,8,7,6,7,6,4,3)
$jValues = @(4,7,5,7,7,6,7,6,6,7,8,7,9,6,7,8,8,7,6,7,7,4,3)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
